{"js": "// \"keep working on discussion\" - results-section wording tweaks.\n//\n//  1. \"...suggesting a shift...\" -> \"...suggesting only a mild shift...\"\n//  2. \"...This suggests that opinions became...\" -> \"...op inions became...\"\n//     (a stray space lands in the middle of \"opinions\")\n//  3. \"...There's a clear increase...\" -> \"...There's a increase...\"\n//     (the word \"clear\" is dropped)\n//  4. A new trailing sentence is appended after \"...on immigration at that time.\"\n//  5. & 6. Two spots where adjacent, identically-formatted runs collapse back\n//     into a single run once the (no-op) text edit is reapplied - no net\n//     text change, just tidies up the run split.\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, replacement) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + searchText);\n  }\n  results.items[0].insertText(replacement, \"Replace\");\n  await context.sync();\n}\n\n// 1) Insert \" only a mild\" before \" shift towards both more extreme ends...\"\nawait replaceOnce(\n  \", suggesting a shift towards both more extreme ends of the opinion spectrum\",\n  \", suggesting only a mild shift towards both more extreme ends of the opinion spectrum\"\n);\n\n// 2) Turn \"opinions\" into \"op inions\" (space inserted mid-word) in the Hungary section.\nawait replaceOnce(\n  \". This suggests that opinions became less neutral and more strongly held in either the positive or negative direction\",\n  \". This suggests that op inions became less neutral and more strongly held in either the positive or negative direction\"\n);\n\n// 3) Drop the word \"clear \" from \"There's a clear increase in non-neutrality...\"\nawait replaceOnce(\n  \": This represents the proportion of individuals holding non-neutral views. There's a clear increase in non-neutrality across all three \",\n  \": This represents the proportion of individuals holding non-neutral views. There's a increase in non-neutrality across all three \"\n);\n\n// 4) Append the new sentence after \"...on immigration at that time.\"\nawait replaceOnce(\n  \"indicating that more people moved away from a neutral stance and adopted a more definite positive or negative opinion on immigration at that time.\",\n  \"indicating that more people moved away from a neutral stance and adopted a more definite positive or negative opinion on immigration at that time. However, shortly after the increases, we see a decrease in non-neutrality until 2020, after which it increases again.\"\n);\n\n// 5) Re-set the first sentence + trailing space of the \"Europe vs. Hungary vs.\n//    Germany\" paragraph so the two runs collapse into one.\nawait replaceOnce(\n  \"we found 2015 to be a turning point in immigration opinion dynamics. \",\n  \"we found 2015 to be a turning point in immigration opinion dynamics. \"\n);\n\n// 6) Re-set \", starting at 2015,\" so the three runs collapse into one.\nawait replaceOnce(\n  \", starting at 2015,\",\n  \", starting at 2015,\"\n);\n\nawait context.sync();\n", "ps1": "# \"keep working on discussion\" - results-section wording tweaks.\n#\n#  1. \"...suggesting a shift...\" -> \"...suggesting only a mild shift...\"\n#  2. \"...This suggests that opinions became...\" -> \"...op inions became...\"\n#     (a stray space lands in the middle of \"opinions\")\n#  3. \"...There's a clear increase...\" -> \"...There's a increase...\"\n#     (the word \"clear\" is dropped)\n#  4. A new trailing sentence is appended after \"...on immigration at that time.\"\n#  5. & 6. Two spots where adjacent, identically-formatted runs collapse back\n#     into a single run once the (no-op) text edit is reapplied - no net\n#     text change, just tidies up the run split.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($findText, $replaceText) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $findText\"\n    }\n}\n\n# 1) Insert \" only a mild\" before \" shift towards both more extreme ends...\"\nReplace-Once \", suggesting a shift towards both more extreme ends of the opinion spectrum\" `\n             \", suggesting only a mild shift towards both more extreme ends of the opinion spectrum\"\n\n# 2) Turn \"opinions\" into \"op inions\" (space inserted mid-word) in the Hungary section.\nReplace-Once \". This suggests that opinions became less neutral and more strongly held in either the positive or negative direction\" `\n             \". This suggests that op inions became less neutral and more strongly held in either the positive or negative direction\"\n\n# 3) Drop the word \"clear \" from \"There's a clear increase in non-neutrality...\"\nReplace-Once \": This represents the proportion of individuals holding non-neutral views. There's a clear increase in non-neutrality across all three \" `\n             \": This represents the proportion of individuals holding non-neutral views. There's a increase in non-neutrality across all three \"\n\n# 4) Append the new sentence after \"...on immigration at that time.\"\nReplace-Once \"indicating that more people moved away from a neutral stance and adopted a more definite positive or negative opinion on immigration at that time.\" `\n             \"indicating that more people moved away from a neutral stance and adopted a more definite positive or negative opinion on immigration at that time. However, shortly after the increases, we see a decrease in non-neutrality until 2020, after which it increases again.\"\n\n# 5) Re-set the first sentence + trailing space of the \"Europe vs. Hungary vs.\n#    Germany\" paragraph so the two runs collapse into one.\nReplace-Once \"we found 2015 to be a turning point in immigration opinion dynamics. \" `\n             \"we found 2015 to be a turning point in immigration opinion dynamics. \"\n\n# 6) Re-set \", starting at 2015,\" so the three runs collapse into one.\nReplace-Once \", starting at 2015,\" `\n             \", starting at 2015,\"\n"}
